$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 60.25
$ws.Range("I8").Value = 60.25
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 180.75
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -41.75
$ws.Range("N8").ClearContents()

$ws.Range("H12").Value = 142.07143
$ws.Range("I12").Value = 138.66667
$ws.Range("J12").Value = 162.5
$ws.Range("K12").Value = 138.66667
$ws.Range("L12").Value = 162.5
$ws.Range("M12").Value = 31.33332999999999
$ws.Range("N12").Value = -502.5

$ws.Range("H33").Value = 266.7931
$ws.Range("I33").Value = 127.478264
$ws.Range("K33").Value = 127.478264
$ws.Range("M33").Value = 101.521736

$ws.Range("H40").Value = 3589.5557
$ws.Range("I40").Value = 4795.4287
$ws.Range("J40").Value = 2290.923
$ws.Range("K40").Value = 4795.4287
$ws.Range("L40").Value = 2290.923
$ws.Range("M40").Value = -4620.4287
$ws.Range("N40").Value = -2640.923

$ws.Range("H132").Value = 89021.14999999999
$ws.Range("I132").Value = 99073.664
$ws.Range("K132").Value = 297220.992
$ws.Range("M132").Value = -294690.992

$ws.Range("H135").Value = 1089.1111
$ws.Range("I135").Value = 889.2692
$ws.Range("K135").Value = 8003.422799999999
$ws.Range("M135").Value = -5468.422799999999

$ws.Range("H137").Value = 857183
$ws.Range("J137").Value = 2032139.5
$ws.Range("L137").Value = 6096418.5
$ws.Range("N137").Value = -6101518.5

$ws.Range("H138").Value = 2783.2083
$ws.Range("J138").Value = 3249.9333
$ws.Range("L138").Value = 9749.7999
$ws.Range("N138").Value = -20029.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8837.895
$ws.Range("I2").Value = 1664.8
$ws.Range("K2").Value = 1664.8
$ws.Range("M2").Value = -1551.8

$ws.Range("H45").Value = 2563.375
$ws.Range("I45").Value = 2563.375
$ws.Range("K45").Value = 2563.375
$ws.Range("M45").Value = -2186.375

$ws.Range("H61").Value = 3031415.2
$ws.Range("I61").Value = 3031415.2
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3031415.2
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3031203.2
$ws.Range("N61").ClearContents()

$ws.Range("H116").Value = 8837.895
$ws.Range("I116").Value = 1664.8
$ws.Range("K116").Value = 1664.8
$ws.Range("M116").Value = 629.2

$ws.Range("H136").Value = 3031415.2
$ws.Range("I136").Value = 3031415.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9094245.600000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9091695.600000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8837.895
$ws.Range("I3").Value = 1664.8
$ws.Range("K3").Value = 1664.8
$ws.Range("M3").Value = -1550.8

$ws.Range("H86").Value = 1649.4286
$ws.Range("I86").Value = 1691.1538
$ws.Range("J86").Value = 1107
$ws.Range("K86").Value = 1691.1538
$ws.Range("L86").Value = 1107
$ws.Range("M86").Value = -568.1538
$ws.Range("N86").Value = -3353

$ws.Range("H89").Value = 1649.4286
$ws.Range("I89").Value = 1691.1538
$ws.Range("J89").Value = 1107
$ws.Range("K89").Value = 8455.769
$ws.Range("L89").Value = 5535
$ws.Range("M89").Value = -2839.769
$ws.Range("N89").Value = -16767

$ws.Range("H105").Value = 4888.6924
$ws.Range("J105").Value = 4000
$ws.Range("L105").Value = 4000
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9687.437
$ws.Range("I31").Value = 3171.8333
$ws.Range("K31").Value = 3171.8333
$ws.Range("M31").Value = -2876.8333

$ws.Range("H34").Value = 9687.437
$ws.Range("I34").Value = 3171.8333
$ws.Range("K34").Value = 3171.8333
$ws.Range("M34").Value = -2969.8333

$ws.Range("H86").Value = 80140.59
$ws.Range("I86").Value = 5737.9287
$ws.Range("J86").Value = 160266.53
$ws.Range("K86").Value = 5737.9287
$ws.Range("L86").Value = 160266.53
$ws.Range("M86").Value = -4614.9287
$ws.Range("N86").Value = -162512.53

$ws.Range("H89").Value = 80140.59
$ws.Range("I89").Value = 5737.9287
$ws.Range("J89").Value = 160266.53
$ws.Range("K89").Value = 28689.6435
$ws.Range("L89").Value = 801332.65
$ws.Range("M89").Value = -23073.6435
$ws.Range("N89").Value = -812564.65

$ws.Range("H132").Value = 10325.579
$ws.Range("I132").Value = 2361.8572
$ws.Range("J132").Value = 32624
$ws.Range("K132").Value = 7085.571599999999
$ws.Range("L132").Value = 97872
$ws.Range("M132").Value = -4555.571599999999
$ws.Range("N132").Value = -102932

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 104011370
$ws.Range("I4").Value = 116887784
$ws.Range("J4").Value = 999999
$ws.Range("K4").Value = 350663352
$ws.Range("L4").Value = 2999997
$ws.Range("M4").Value = -350663240
$ws.Range("N4").Value = -3000221

$ws.Range("H40").Value = 119.37037
$ws.Range("J40").Value = 146.5625
$ws.Range("L40").Value = 586.25
$ws.Range("N40").Value = -724.25

$ws.Range("H129").Value = 1673.5555
$ws.Range("J129").Value = 3994.5
$ws.Range("L129").Value = 11983.5
$ws.Range("N129").Value = -21983.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 30472.436
$ws.Range("I122").Value = 41124.73
$ws.Range("K122").Value = 123374.19
$ws.Range("M122").Value = -120924.19

$ws.Range("H132").Value = 933947.5600000001
$ws.Range("I132").Value = 1101319.9
$ws.Range("J132").Value = 13400
$ws.Range("K132").Value = 3303959.7
$ws.Range("L132").Value = 40200
$ws.Range("M132").Value = -3301429.7
$ws.Range("N132").Value = -45260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3890.4583
$ws.Range("I7").Value = 3609.8823
$ws.Range("K7").Value = 3609.8823
$ws.Range("M7").Value = -3497.8823

$ws.Range("H40").Value = 4047.2273
$ws.Range("I40").Value = 4155.4707
$ws.Range("K40").Value = 4155.4707
$ws.Range("M40").Value = -4019.4707

$ws.Range("H68").Value = 3542.5
$ws.Range("J68").Value = 3145.1428
$ws.Range("L68").Value = 3145.1428
$ws.Range("N68").Value = -4643.1428

$ws.Range("H71").Value = 3542.5
$ws.Range("J71").Value = 3145.1428
$ws.Range("L71").Value = 15725.714
$ws.Range("N71").Value = -23213.714

$ws.Range("H126").Value = 3890.4583
$ws.Range("I126").Value = 3609.8823
$ws.Range("K126").Value = 10829.6469
$ws.Range("M126").Value = -8359.6469

$ws.Range("H132").Value = 1663042.4
$ws.Range("I132").Value = 2179899.5
$ws.Range("J132").Value = 9099.799999999999
$ws.Range("K132").Value = 6539698.5
$ws.Range("L132").Value = 27299.4
$ws.Range("M132").Value = -6537168.5
$ws.Range("N132").Value = -32359.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H126").Value = 4274.5
$ws.Range("I126").Value = 3979.45
$ws.Range("J126").Value = 5749.75
$ws.Range("K126").Value = 11938.35
$ws.Range("L126").Value = 17249.25
$ws.Range("M126").Value = -9468.349999999999
$ws.Range("N126").Value = -22189.25

$ws.Range("H132").Value = 7189003.5
$ws.Range("I132").Value = 7741542
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 23224626
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -23222096
$ws.Range("N132").Value = -23060
